# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve profit-tracking tables across all eight
# Disciple-of-the-Hand job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected row the currentAveragePrice / currentAveragePriceNQ / 
# currentAveragePriceHQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# values (columns H:N) are refreshed with newly observed market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 12 (Leve Item ID 5515)
$ws.Range("H12").Value2 = 165.2
$ws.Range("I12").Value2 = 154
$ws.Range("J12").Value2 = 210
$ws.Range("K12").Value2 = 154
$ws.Range("L12").Value2 = 210
$ws.Range("M12").Value2 = 16
$ws.Range("N12").Value2 = -550

# ALC row 53 (Leve Item ID 5479)
$ws.Range("H53").Value2 = 180.68182
$ws.Range("I53").Value2 = 191.46153
$ws.Range("K53").Value2 = 191.46153
$ws.Range("M53").Value2 = 445.53847

# ALC row 70 (Leve Item ID 12604)
$ws.Range("H70").Value2 = 3695.5
$ws.Range("I70").Value2 = 2751
$ws.Range("J70").Value2 = 4640
$ws.Range("K70").Value2 = 8253
$ws.Range("L70").Value2 = 13920
$ws.Range("M70").Value2 = -7983
$ws.Range("N70").Value2 = -14460

# ALC row 73 (Leve Item ID 12604)
$ws.Range("H73").Value2 = 3695.5
$ws.Range("I73").Value2 = 2751
$ws.Range("J73").Value2 = 4640
$ws.Range("K73").Value2 = 8253
$ws.Range("L73").Value2 = 13920
$ws.Range("M73").Value2 = -7317
$ws.Range("N73").Value2 = -15792

# ALC row 88 (Leve Item ID 12608)
$ws.Range("H88").Value2 = 3547.762
$ws.Range("I88").Value2 = 4574.875
$ws.Range("J88").Value2 = 2915.6924
$ws.Range("K88").Value2 = 4574.875
$ws.Range("L88").Value2 = 2915.6924
$ws.Range("M88").Value2 = -4168.875
$ws.Range("N88").Value2 = -3727.6924

# ALC row 91 (Leve Item ID 12608)
$ws.Range("H91").Value2 = 3547.762
$ws.Range("I91").Value2 = 4574.875
$ws.Range("J91").Value2 = 2915.6924
$ws.Range("K91").Value2 = 4574.875
$ws.Range("L91").Value2 = 2915.6924
$ws.Range("M91").Value2 = -3170.875
$ws.Range("N91").Value2 = -5723.6924

# ALC row 125 (Leve Item ID 36228)
$ws.Range("H125").Value2 = 2333.3333
$ws.Range("J125").Value2 = 0
$ws.Range("L125").Value2 = 0
$ws.Range("N125").ClearContents()

# ALC row 132 (Leve Item ID 44049)
$ws.Range("H132").Value2 = 5331.174
$ws.Range("I132").Value2 = 5854.5854
$ws.Range("K132").Value2 = 17563.7562
$ws.Range("M132").Value2 = -15033.7562

# ALC row 135 (Leve Item ID 44047)
$ws.Range("H135").Value2 = 5517.857
$ws.Range("I135").Value2 = 276.75
$ws.Range("K135").Value2 = 2490.75
$ws.Range("M135").Value2 = 44.25

# ALC row 138 (Leve Item ID 44169)
$ws.Range("H138").Value2 = 3217.1128
$ws.Range("I138").Value2 = 4422.5713
$ws.Range("J138").Value2 = 2710.82
$ws.Range("K138").Value2 = 13267.7139
$ws.Range("L138").Value2 = 8132.460000000001
$ws.Range("M138").Value2 = -8127.713899999999
$ws.Range("N138").Value2 = -18412.46

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61 (Leve Item ID 43999)
$ws.Range("H61").Value2 = 2706031.8
$ws.Range("I61").Value2 = 3198.4243
$ws.Range("K61").Value2 = 3198.4243
$ws.Range("M61").Value2 = -2986.4243

# ARM row 108 (Leve Item ID 27084)
$ws.Range("H108").Value2 = 49999.5
$ws.Range("I108").Value2 = 0
$ws.Range("J108").Value2 = 49999.5
$ws.Range("K108").Value2 = 0
$ws.Range("L108").Value2 = 49999.5
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value2 = -57679.5

# ARM row 136 (Leve Item ID 43999)
$ws.Range("H136").Value2 = 2706031.8
$ws.Range("I136").Value2 = 3198.4243
$ws.Range("K136").Value2 = 9595.2729
$ws.Range("M136").Value2 = -7045.2729

$ws = $wb.Worksheets.Item("BSM")
# BSM row 99 (Leve Item ID 19943)
$ws.Range("H99").Value2 = 5635.769
$ws.Range("I99").Value2 = 7818.625
$ws.Range("J99").Value2 = 2143.2
$ws.Range("K99").Value2 = 7818.625
$ws.Range("L99").Value2 = 2143.2
$ws.Range("M99").Value2 = -6320.625
$ws.Range("N99").Value2 = -5139.2

# BSM row 134 (Leve Item ID 43998)
$ws.Range("H134").Value2 = 4764096.5
$ws.Range("I134").Value2 = 2340.3438
$ws.Range("K134").Value2 = 7021.0314
$ws.Range("M134").Value2 = -4486.0314

$ws = $wb.Worksheets.Item("CRP")
# CRP row 11 (Leve Item ID 1821)
$ws.Range("H11").Value2 = 1995
$ws.Range("J11").Value2 = 1995
$ws.Range("L11").Value2 = 1995
$ws.Range("N11").Value2 = -2275

# CRP row 12 (Leve Item ID 1604)
$ws.Range("H12").Value2 = 2843
$ws.Range("I12").Value2 = 186
$ws.Range("J12").Value2 = 5500
$ws.Range("K12").Value2 = 186
$ws.Range("L12").Value2 = 5500
$ws.Range("M12").Value2 = -16
$ws.Range("N12").Value2 = -5840

# CRP row 99 (Leve Item ID 36198)
$ws.Range("H99").Value2 = 8598.471
$ws.Range("I99").Value2 = 18242.154
$ws.Range("J99").Value2 = 2628.5715
$ws.Range("K99").Value2 = 18242.154
$ws.Range("L99").Value2 = 2628.5715
$ws.Range("M99").Value2 = -16744.154
$ws.Range("N99").Value2 = -5624.5715

# CRP row 126 (Leve Item ID 36198)
$ws.Range("H126").Value2 = 8598.471
$ws.Range("I126").Value2 = 18242.154
$ws.Range("J126").Value2 = 2628.5715
$ws.Range("K126").Value2 = 54726.462
$ws.Range("L126").Value2 = 7885.7145
$ws.Range("M126").Value2 = -52256.462
$ws.Range("N126").Value2 = -12825.7145

# CRP row 132 (Leve Item ID 44019)
$ws.Range("H132").Value2 = 3545.037
$ws.Range("I132").Value2 = 3424.3809
$ws.Range("J132").Value2 = 3967.3333
$ws.Range("K132").Value2 = 10273.1427
$ws.Range("L132").Value2 = 11901.9999
$ws.Range("M132").Value2 = -7743.1427
$ws.Range("N132").Value2 = -16961.9999

# CRP row 134 (Leve Item ID 44020)
$ws.Range("H134").Value2 = 1720.2963
$ws.Range("I134").Value2 = 1638.875
$ws.Range("K134").Value2 = 4916.625
$ws.Range("M134").Value2 = -2381.625

$ws = $wb.Worksheets.Item("CUL")
# CUL row 8 (Leve Item ID 16734)
$ws.Range("H8").Value2 = 288.2
$ws.Range("I8").Value2 = 288.2
$ws.Range("K8").Value2 = 864.5999999999999
$ws.Range("M8").Value2 = -725.5999999999999

# CUL row 12 (Leve Item ID 4854)
$ws.Range("H12").Value2 = 957.38464
$ws.Range("J12").Value2 = 774.375
$ws.Range("L12").Value2 = 2323.125
$ws.Range("N12").Value2 = -2669.125

# CUL row 82 (Leve Item ID 12856)
$ws.Range("H82").Value2 = 13929.75
$ws.Range("I82").Value2 = 10013
$ws.Range("J82").Value2 = 14489.286
$ws.Range("K82").Value2 = 30039
$ws.Range("L82").Value2 = 43467.858
$ws.Range("M82").Value2 = -29633
$ws.Range("N82").Value2 = -44279.858

# CUL row 85 (Leve Item ID 12856)
$ws.Range("H85").Value2 = 13929.75
$ws.Range("I85").Value2 = 10013
$ws.Range("J85").Value2 = 14489.286
$ws.Range("K85").Value2 = 30039
$ws.Range("L85").Value2 = 43467.858
$ws.Range("M85").Value2 = -28635
$ws.Range("N85").Value2 = -46275.858

# CUL row 92 (Leve Item ID 19841)
$ws.Range("H92").Value2 = 765.2593000000001
$ws.Range("I92").Value2 = 453.1
$ws.Range("K92").Value2 = 1359.3
$ws.Range("M92").Value2 = -111.3000000000002

# CUL row 110 (Leve Item ID 27857)
$ws.Range("H110").Value2 = 18666.666
$ws.Range("I110").Value2 = 12000
$ws.Range("K110").Value2 = 36000
$ws.Range("M110").Value2 = -31910

# CUL row 120 (Leve Item ID 27877)
$ws.Range("H120").Value2 = 20456.818
$ws.Range("I120").Value2 = 14170.833
$ws.Range("J120").Value2 = 28000
$ws.Range("K120").Value2 = 42512.499
$ws.Range("L120").Value2 = 84000
$ws.Range("M120").Value2 = -37674.499
$ws.Range("N120").Value2 = -93676

# CUL row 131 (Leve Item ID 36060)
$ws.Range("H131").Value2 = 3870.9048
$ws.Range("I131").Value2 = 942.2857
$ws.Range("J131").Value2 = 5335.2144
$ws.Range("K131").Value2 = 2826.8571
$ws.Range("L131").Value2 = 16005.6432
$ws.Range("M131").Value2 = 2213.1429
$ws.Range("N131").Value2 = -26085.6432

# CUL row 133 (Leve Item ID 44073)
$ws.Range("H133").Value2 = 6140.9414
$ws.Range("I133").Value2 = 3722.7693
$ws.Range("J133").Value2 = 14000
$ws.Range("K133").Value2 = 11168.3079
$ws.Range("L133").Value2 = 42000
$ws.Range("M133").Value2 = -6108.3079
$ws.Range("N133").Value2 = -52120

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70 (Leve Item ID 14146)
$ws.Range("H70").Value2 = 34441
$ws.Range("I70").Value2 = 33604.125
$ws.Range("K70").Value2 = 33604.125
$ws.Range("M70").Value2 = -33334.125

# GSM row 73 (Leve Item ID 14146)
$ws.Range("H73").Value2 = 34441
$ws.Range("I73").Value2 = 33604.125
$ws.Range("K73").Value2 = 33604.125
$ws.Range("M73").Value2 = -32668.125

# GSM row 80 (Leve Item ID 12521)
$ws.Range("H80").Value2 = 1896.1
$ws.Range("I80").Value2 = 1702
$ws.Range("J80").Value2 = 2090.2
$ws.Range("K80").Value2 = 1702
$ws.Range("L80").Value2 = 2090.2
$ws.Range("M80").Value2 = -704
$ws.Range("N80").Value2 = -4086.2

# GSM row 83 (Leve Item ID 12521)
$ws.Range("H83").Value2 = 1896.1
$ws.Range("I83").Value2 = 1702
$ws.Range("J83").Value2 = 2090.2
$ws.Range("K83").Value2 = 8510
$ws.Range("L83").Value2 = 10451
$ws.Range("M83").Value2 = -3518
$ws.Range("N83").Value2 = -20435

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16 (Leve Item ID 5289)
$ws.Range("H16").Value2 = 397.55
$ws.Range("I16").Value2 = 413.27777
$ws.Range("K16").Value2 = 413.27777
$ws.Range("M16").Value2 = -243.27777

# LTW row 34 (Leve Item ID 3347)
$ws.Range("H34").Value2 = 10082.857
$ws.Range("J34").Value2 = 13660
$ws.Range("L34").Value2 = 13660
$ws.Range("N34").Value2 = -14004

# LTW row 56 (Leve Item ID 3668)
$ws.Range("H56").Value2 = 11000
$ws.Range("I56").Value2 = 10000
$ws.Range("J56").Value2 = 12000
$ws.Range("K56").Value2 = 10000
$ws.Range("L56").Value2 = 12000
$ws.Range("M56").Value2 = -9309
$ws.Range("N56").Value2 = -13382

$ws = $wb.Worksheets.Item("WVR")
# WVR row 81 (Leve Item ID 12596)
$ws.Range("H81").Value2 = 93809.27
$ws.Range("I81").Value2 = 2783.3333
$ws.Range("J81").Value2 = 203040.4
$ws.Range("K81").Value2 = 5566.6666
$ws.Range("L81").Value2 = 406080.8
$ws.Range("M81").Value2 = -4505.6666
$ws.Range("N81").Value2 = -408202.8

# WVR row 84 (Leve Item ID 12596)
$ws.Range("H84").Value2 = 93809.27
$ws.Range("I84").Value2 = 2783.3333
$ws.Range("J84").Value2 = 203040.4
$ws.Range("K84").Value2 = 27833.333
$ws.Range("L84").Value2 = 2030404
$ws.Range("M84").Value2 = -22529.333
$ws.Range("N84").Value2 = -2041012

# WVR row 126 (Leve Item ID 36210)
$ws.Range("H126").Value2 = 3550
$ws.Range("J126").Value2 = 3550
$ws.Range("L126").Value2 = 10650
$ws.Range("N126").Value2 = -15590

# WVR row 132 (Leve Item ID 44029)
$ws.Range("H132").Value2 = 6411846
$ws.Range("I132").Value2 = 7577136.5
$ws.Range("K132").Value2 = 22731409.5
$ws.Range("M132").Value2 = -22728879.5

# WVR row 136 (Leve Item ID 44031)
$ws.Range("H136").Value2 = 4609862.5
$ws.Range("I136").Value2 = 2900264.8
$ws.Range("J136").Value2 = 7815359
$ws.Range("K136").Value2 = 8700794.399999999
$ws.Range("L136").Value2 = 23446077
$ws.Range("M136").Value2 = -8698244.399999999
$ws.Range("N136").Value2 = -23451177

